$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new (blank) column at G - shifts old G,H,I to H,I,J
$ws.Columns("G").Insert()

# The insert operation clones the header cell into the new gap (G1);
# the source table doesn't have a header there, so drop it entirely.
$ws.Range("G1").Clear()

# 2. Widen column B to fit the longer labels
$ws.Columns("B").ColumnWidth = 33

# 3. Append the six new "elderly" breakdown rows, copying the formatting
#    from the last existing data row (16) first so styles/number formats match.
$newRows = @(
    @{ Row=17; A=15; B="admissoes_gerais_non_elderly";          C=5143596; D=4527062; E=4763975; F=-7.38045911848442;  H=-3.86908379352856;  I=-3.92982532094831;  J=-3.808303861545 },
    @{ Row=18; A=16; B="admissoes_gerais_uti_non_elderly";      C=182960;  D=211836;  E=313738;  F=71.4790118058592;   H=32.3763878656542;   I=31.9920084301855;   J=32.7618866685162 },
    @{ Row=19; A=17; B="admissoes_gerais_non_uti_non_elderly";  C=4960636; D=4315226; E=4450237; F=-10.2889831061985;  H=-5.42772146361673;  I=-5.48902439280426;  J=-5.3663787713318 },
    @{ Row=20; A=18; B="admissoes_gerais_elderly";              C=2525487; D=2228491; E=2425462; F=-3.960622248303;    H=-2.06828530773551;  I=-2.15599037230442;  J=-1.98050162641347 },
    @{ Row=21; A=19; B="admissoes_gerais_uti_elderly";          C=269248;  D=326352;  E=388800;  F=44.4021868314714;   H=20.1038987707989;   I=19.8112374721649;   J=20.3972749492516 },
    @{ Row=22; A=20; B="admissoes_gerais_non_uti_elderly";      C=2256239; D=1902139; E=2036662; F=-9.73199204516897;  H=-5.18011746468455;  I=-5.27158534413688;  J=-5.08856126566094 }
)

$ws.Range("A16:J16").Copy()

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row`:J$row").PasteSpecial(-4122)

    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = $r.F
    $ws.Range("H$row").Value = $r.H
    $ws.Range("I$row").Value = $r.I
    $ws.Range("J$row").Value = $r.J
}

# 4. Match the author's final selection
$ws.Range("R8").Select()
